# Applies the text edits described by the diff to the promotion-exam
# presentation. The presentation in this environment only contains the
# four "content" slides (title, "leader", "subcontract", "grow with the
# company"); the many datetimeFigureOut field updates (2020/2/14 ->
# 2020/2/17) in the original diff all live on slides that are not part
# of this deck, so there is nothing to change for those hunks here.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 2 - "リーダーとしての活躍"
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$body2 = $s2.Shapes.Item(2).TextFrame.TextRange

# "プロジェクトのリーダー" + "として" + "、" -> single run
# "プロジェクトのリーダーとして、" (15 characters from the start of the body)
$run = $body2.Characters(1, 15)
$run.Text = "プロジェクトのリーダーとして、"

# "・技術の展開" -> "・" + "生産性向上と品質向上のため、技術を検討し、メンバに展開"
# Locate it via the 4th paragraph, which reads "　・技術の展開".
$para4 = $body2.Paragraphs(4)
$tail = $para4.Characters(3, 5)
$tail.Text = "生産性向上と品質向上のため、技術を検討し、メンバに展開"

# Last paragraph "　・" -> "　"
$body2b = $s2.Shapes.Item(2).TextFrame.TextRange
$lastPara = $body2b.Paragraphs($body2b.Paragraphs().Count)
$lastPara.Text = "　"

# ---------------------------------------------------------------------
# Slide 3 - "請負の活躍"
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$body3 = $s3.Shapes.Item(2).TextFrame.TextRange

# "本業以外に、請負社内開発も兼務した。" -> four runs:
#   "本業以外に、請負社内" / "開発" / "の" / "兼務。"
$para1 = $body3.Paragraphs(1)
$r1 = $para1.Characters(1, 10)
$r1.Text = "本業以外に、請負社内"
$r2 = $para1.Characters(11, 2)
$r2.Text = "開発"
$r3 = $para1.Characters(13, 1)
$r3.Text = "の"
$r4 = $para1.Characters(14, 5)
$r4.Text = "兼務。"

# "・顧客と仕様の確認" -> "・顧客と仕様" + "の調整"
$body3b = $s3.Shapes.Item(2).TextFrame.TextRange
$para2 = $body3b.Paragraphs(2)
$r5 = $para2.Characters(2, 6)
$r5.Text = "・顧客と仕様"
$r6 = $para2.Characters(8, 3)
$r6.Text = "の調整"

# ---------------------------------------------------------------------
# Slide 4 - "会社と共に成長"
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$body4 = $s4.Shapes.Item(2).TextFrame.TextRange

# "マネジメント、" + "請負経験を生かして、" -> single run
$run4 = $body4.Characters(1, 17)
$run4.Text = "マネジメント、請負経験を生かして、"
